# 57817 - Implement Dominium employee logic in the EDC, consistent with the
# logic already implemented in ICW (see Findings row 17: ICW_DominiumEmployee).

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsFindings = $wb.Worksheets.Item("Findings")

# --- Findings sheet: add the new EDC Dominium-employee finding row -----
# Mirrors the existing ICW_DominiumEmployee row (row 17): key / message /
# "Newly added by Raluca" marker. Insert a fresh row 137 (shifting the
# following rows down by one) and populate it.
$wsFindings.Range("A137").EntireRow.Insert()
$wsFindings.Range("A137").Value = "EDC_DominiumEmployee"
$wsFindings.Range("B137").Value = "Applicant is a Dominium employee, manual review is required."
$wsFindings.Range("C137").Value = "Newly added by Raluca"

# --- Settings sheet: refresh the local prompts folder path for the new --
# --- author's machine ----------------------------------------------------
$wsSettings.Range("B7").Value = "C:\Users\raluca.ilinca.AzureAI-Jump2\Documents\AZApplicationReview\Prompts\"

# --- View/selection state as left by the author ------------------------
$wsFindings.Select()
$excel.ActiveWindow.ScrollRow = 114
$excel.ActiveWindow.ScrollColumn = 1
$wsFindings.Range("D137").Select()

$wsSettings.Activate()
$wsSettings.Range("B7").Select()
